# Applies the cryptos-list refresh described by the commit diff:
# updates the Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "30.237.98"; ForceText = $false },
    @{ Cell = "D3"; Value = "1.985.31"; ForceText = $false },
    @{ Cell = "E3"; Value = "  +5.84%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  -0.20%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "323.36"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +1.08%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "1.000"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -0.15%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.5100"; ForceText = $true },
    @{ Cell = "E7"; Value = "  +1.06%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "0.4104"; ForceText = $true },
    @{ Cell = "E8"; Value = "  +3.60%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.08850"; ForceText = $true },
    @{ Cell = "E9"; Value = "  +7.80%  "; ForceText = $false },
    @{ Cell = "E10"; Value = "  +3.21%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "42.44"; ForceText = $true },
    @{ Cell = "E11"; Value = "  +0.78%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "24.17"; ForceText = $true },
    @{ Cell = "E12"; Value = "  +2.58%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "1.981.48"; ForceText = $false },
    @{ Cell = "E13"; Value = "  +6.05%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "6.485"; ForceText = $true },
    @{ Cell = "E14"; Value = "  +3.09%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "7.387"; ForceText = $true },
    @{ Cell = "E15"; Value = "  +2.69%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "1.001"; ForceText = $true },
    @{ Cell = "E16"; Value = "  -0.08%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "93.90"; ForceText = $true },
    @{ Cell = "E17"; Value = "  +2.13%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "0.00001122"; ForceText = $true },
    @{ Cell = "E18"; Value = "  +3.32%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "0.06542"; ForceText = $true },
    @{ Cell = "E19"; Value = "  +0.93%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "18.79"; ForceText = $true },
    @{ Cell = "E20"; Value = "  +3.67%  "; ForceText = $false },
    @{ Cell = "E21"; Value = "  -0.09%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "6.073"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +4.19%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "30.293.63"; ForceText = $false },
    @{ Cell = "E23"; Value = "  +0.75%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "11.47"; ForceText = $true },
    @{ Cell = "E24"; Value = "  +3.06%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "2.212"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +2.04%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "2.211.65"; ForceText = $false },
    @{ Cell = "E26"; Value = "  +6.07%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "22.43"; ForceText = $true },
    @{ Cell = "E27"; Value = "  +5.90%  "; ForceText = $false },
    @{ Cell = "E28"; Value = "  +1.40%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "2.363"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +5.60%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "130.69"; ForceText = $true },
    @{ Cell = "E30"; Value = "  +2.62%  "; ForceText = $false },
    @{ Cell = "E31"; Value = "  +4.74%  "; ForceText = $false },
    @{ Cell = "E32"; Value = "  +1.81%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "6.035"; ForceText = $true },
    @{ Cell = "E33"; Value = "  +1.57%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "3.806"; ForceText = $true },
    @{ Cell = "E34"; Value = "  +2.77%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "1.312"; ForceText = $true },
    @{ Cell = "D36"; Value = "0.02480"; ForceText = $true },
    @{ Cell = "E36"; Value = "  +2.05%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "5.376"; ForceText = $true },
    @{ Cell = "E37"; Value = "  +1.87%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "0.06490"; ForceText = $true },
    @{ Cell = "E38"; Value = "  +1.86%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "0.2175"; ForceText = $true },
    @{ Cell = "E39"; Value = "  +1.85%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "8.902"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +4.64%  "; ForceText = $false },
    @{ Cell = "E41"; Value = "  +4.23%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "11.75"; ForceText = $true },
    @{ Cell = "E42"; Value = "  +4.39%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "1.218"; ForceText = $true },
    @{ Cell = "E43"; Value = "  +0.28%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "13.61"; ForceText = $true },
    @{ Cell = "E44"; Value = "  +3.59%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "0.6118"; ForceText = $true },
    @{ Cell = "E45"; Value = "  +3.57%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "2.186"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +4.50%  "; ForceText = $false },
    @{ Cell = "E47"; Value = "  +0.71%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "123.69"; ForceText = $true },
    @{ Cell = "E48"; Value = "  +1.26%  "; ForceText = $false },
    @{ Cell = "E49"; Value = "  +1.12%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "79.41"; ForceText = $true },
    @{ Cell = "E50"; Value = "  +2.43%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "0.06869"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +2.06%  "; ForceText = $false }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        # A leading apostrophe forces Excel to store an otherwise numeric-looking
        # string (e.g. "1.000") as literal text instead of coercing it to a number.
        $ws.Range($u.Cell).Value = "'" + $u.Value
        # Re-apply the default style so the quote-prefix formatting Excel applies
        # when forcing text does not leave a stray style on the cell.
        $ws.Range($u.Cell).Style = "Normal"
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
